$d = $word.ActiveDocument
$vt = [char]11

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Host "WARNING: not found -> [$find]"
    }
    return $result
}

# --- Title ---
Replace-Text "The Enigma of Consciousness: Exploring Boundless Mysteries" "Exploring the Realm of Biology: Unveiling the Secrets of Life"

# --- Author name (Dr. Nathan Green -> Mary Taylor), merges 3 runs into 1 ---
Replace-Text "Dr. Nathan Green" "Mary Taylor"

# --- Email local-part and domain suffix (keep the "." run untouched) ---
Replace-Text "nathangreen@domainname" "marytaylor@validhayschool"
Replace-Text "org" "edu"

# --- Body paragraph 1 edits ---
Replace-Text "As sentient beings, we find ourselves immersed in the profound sea of consciousness, a phenomenon that has intrigued philosophers, scientists, and artists for centuries" "In the heart of the natural world, where complex systems intertwine and life's mysteries unfold, lies the captivating study of biology"

Replace-Text " Our waking thoughts, dreams, feelings, and perceptions form a tapestry of experience that defines our existence" " From the intricate workings of cells to the grand tapestry of ecosystems, biology unveils the essence of existence, revealing the profound interconnectedness of all living organisms"

$old1 = " What, then, is the nature of consciousness? How do physical processes in our brain give rise to subjective experience?" + $vt + $vt + "In the pursuit of understanding consciousness, we journey through an array of disciplines. Neuroscientists delve into the intricate neural networks of the brain, seeking correlations between brain activity and conscious states. Psychologists examine the role of attention, memory, and other cognitive processes in shaping our perceptions. Philosophers contemplate the hard problem of consciousness, grappling with the fundamental question of existence of subjective experience"
$new1 = " In this journey of discovery, we embark on an exploration of the fundamental principles that govern life, seeking to comprehend the symphony of processes that shape the natural world around us"
Replace-Text $old1 $new1

Replace-Text "From mystics seeking eternal truths to artists exploring the depths of the human condition, consciousness has permeated diverse fields of study" "As we delve into the microscopic realm, we unravel the secrets hidden within the building blocks of life - cells"

Replace-Text " Poets capture the elusive essence of awareness in their verse, while musicians weave melodies that evoke emotions and touch the soul" " These minuscule entities, though unseen to the naked eye, orchestrate a symphony of functions that sustain and govern all living organisms"

$old2 = " Cultural anthropologists probe the variations in consciousness across different societies, revealing the influence of language, ritual, and belief systems"
$new2 = " Through the study of cellular processes, we gain insights into the mechanisms that control growth, reproduction, and inheritance, marveling at the intricate dance of DNA molecules that holds the blueprint for life." + $vt + $vt + "Further, biology unveils the complexities of life's adaptations, showcasing how organisms have evolved over time to survive in diverse environments. We witness the breathtaking strategies employed by creatures great and small, from the intricate camouflage of a chameleon to the sleek design of a dolphin. These adaptations testify to the resilience and ingenuity of life, as organisms constantly strive to harmonize with their surroundings." + $vt + $vt + "In the realm of biology, we encounter the awe-inspiring tapestry of ecosystems, where a multitude of organisms interact in intricate relationships. Food chains and webs delineate the intricate dance of energy and nutrient transfer, demonstrating the interconnectedness of all living things. We explore the delicate balance that sustains these ecosystems, highlighting the profound impact human activities have on the health and stability of our planet"
Replace-Text $old2 $new2

# --- Summary heading: add lastRenderedPageBreak before "Summary" text ---
$d.Paragraphs(6).Range.Find.Execute("Summary", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# --- Summary paragraph edits ---
Replace-Text "The enigma of consciousness remains an alluring and inscrutable aspect of human existence" "Biology unveils the mysteries of life, delving into the intricacies of cells, adaptations, and ecosystems"

Replace-Text " Its multidimensional nature defies easy explanation, inviting exploration from a wide array of disciplines" " The study of biology equips us with a profound understanding of the symphony of life, emphasizing the interconnectedness of all living organisms and the delicate balance of nature"

$old3 = " As we continue to probe the mysteries of consciousness, we gain a deeper appreciation for the complexities of the human mind and the boundless realm of subjective experience. The pursuit of understanding consciousness is not merely an intellectual endeavor; it is an exploration of our very being"
$new3 = " With each discovery, we deepen our appreciation for the beauty and complexity of life, fostering a sense of stewardship and responsibility towards the natural world that sustains us"
Replace-Text $old3 $new3

Write-Host "Done with text replacements"
